# Jake Alinabon Salibay Q0554 - "adding averages and more checks"
#
# 1) Training Dashboard sheet: refresh the "PERIOD TO EXPIRE" (H) and
#    "LAST UPDATE" (I) columns for rows 3-17 (the data moved forward
#    8 days, from 08-Sep-2025 to 16-Sep-2025).
# 2) Header/title formatting: header row text turns white (readable on
#    the dark-blue fill) and the title no longer uses an oversized font.
# 3) Exam Dashboard sheet: widen the COMMENTS-adjacent category column and
#    reword the remark in E3.

$wb = $excel.ActiveWorkbook

$training = $wb.Worksheets.Item("Training Dashboard")
$exam     = $wb.Worksheets.Item("Exam Dashboard")

# --- Training Dashboard: updated "period to expire" + "last update" ----
$periodToExpire = @{
    3  = 241
    4  = 231
    5  = 310
    6  = 334
    7  = 316
    8  = 316
    9  = 316
    10 = 247
    11 = 266
    12 = 245
    13 = 210
    14 = 316
    15 = 317
    16 = 316
    17 = -180
}

foreach ($row in 3..17) {
    $training.Cells.Item($row, 8).Value = $periodToExpire[$row]
    # Leading apostrophe keeps this a literal text value (matches the
    # original inline-string "dd-MMM-yyyy" cells) instead of Excel
    # auto-converting the date-shaped text into a real date serial.
    $training.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# --- Header styling: bold white text for the dark header row, and drop
#     the oversized title font so it matches the header's weight -------
$training.Range("A2:K2").Font.Bold = $true
$training.Range("A2:K2").Font.Color = 16777215
$training.Range("A1").Font.Size = 11
$training.Range("A1").Font.Color = 16777215

# --- Exam Dashboard: widen column E and update the remark text --------
$exam.Columns.Item(5).ColumnWidth = 14.17
$exam.Range("E3").Value = "date is valid"
